$wb = $excel.ActiveWorkbook

# Remember which sheet was active so we can restore the selection at the end.
$originalActive = $wb.ActiveSheet.Name

# --- Insert the new "2022-Q1" worksheet right before the "总计" sheet ---
$total = $wb.Worksheets.Item("总计")
$template = $wb.Worksheets.Item("2021-Q4")

$newSheet = $wb.Worksheets.Add($total)
$newSheet.Name = "2022-Q1"

# Match the look & feel (fonts/borders/margins/outline) of the sibling
# quarterly sheets instead of Excel's blank-sheet defaults.
$newSheet.Outline.SummaryRow = 1
$newSheet.Outline.SummaryColumn = 1
$newSheet.PageSetup.LeftMargin = $template.PageSetup.LeftMargin
$newSheet.PageSetup.RightMargin = $template.PageSetup.RightMargin
$newSheet.PageSetup.TopMargin = $template.PageSetup.TopMargin
$newSheet.PageSetup.BottomMargin = $template.PageSetup.BottomMargin
$newSheet.PageSetup.HeaderMargin = $template.PageSetup.HeaderMargin
$newSheet.PageSetup.FooterMargin = $template.PageSetup.FooterMargin

$template.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Header row (same layout as the other quarterly sheets)
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Force the data columns that must stay textual (codes & percentages kept
# as strings in the source data) to Text format before writing them, so
# Excel doesn't silently coerce them into numbers; then drop back to the
# plain "Normal" style so no stray formatting is left behind on the cells.
$newSheet.Range("B2:G3").NumberFormat = "@"

# Row 2 - 新疆前海联合泳涛灵活配置混合A
$newSheet.Range("B2").Value = "004634"
$newSheet.Range("C2").Value = "新疆前海联合泳涛灵活配置混合A"
$newSheet.Range("D2").Value = "1.33"
$newSheet.Range("E2").Value = "89.65"
$newSheet.Range("F2").Value = "4.38"
$newSheet.Range("G2").Value = "0.0583"

# Row 3 - 新疆前海联合泳涛灵活配置混合C
$newSheet.Range("B3").Value = "007041"
$newSheet.Range("C3").Value = "新疆前海联合泳涛灵活配置混合C"
$newSheet.Range("D3").Value = "0.00"
$newSheet.Range("E3").Value = "89.65"
$newSheet.Range("F3").Value = "4.38"

$newSheet.Range("B2:G3").Style = "Normal"

$newSheet.Range("A2").Value = 0
$newSheet.Range("H2").Value = 7
$newSheet.Range("A3").Value = 1
$newSheet.Range("G3").Value = 0
$newSheet.Range("H3").Value = 7

# --- Update the "总计" summary sheet: insert a new top row for 2022-Q1 ---
$summary = $wb.Worksheets.Item("总计")
$summary.Rows.Item(2).Insert()

# The inserted row inherits the header row's style on B:D; strip that back
# to the plain (unstyled) look the data rows use.
$summary.Range("B2:D2").ClearFormats()

# Column A carries the bold/centered "index" style on every data row; copy
# it from a still-intact sibling cell onto the freshly inserted A2.
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q1"
$summary.Range("C2").Value = 2
$summary.Range("D2").Value = 0.06

# Renumber the sequence index in column A for the rows that shifted down.
$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2
$summary.Range("A5").Value = 3
$summary.Range("A6").Value = 4

# Restore the originally active sheet/selection.
$wb.Worksheets.Item($originalActive).Activate()
